$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Product Backlog (Estimativa): fill in the TAMANHO (E) and TAMANHO(NUMBER) (F)
# columns for each requirement row in the "Requisitos" table.

# Row 2 - Tela Cadastro
$ws.Range("E2").Value = "G"
$ws.Range("F2").Value = 13

# Row 3 - Tela Login
$ws.Range("E3").Value = "M"
$ws.Range("F3").Value = 8

# Row 4 - Cabeçalho
$ws.Range("E4").Value = "P"
$ws.Range("F4").Value = 5

# Row 5 - Tela Inicial
$ws.Range("E5").Value = "G"
$ws.Range("F5").Value = 13

# Row 6 - Footer
$ws.Range("E6").Value = "M"
$ws.Range("F6").Value = 8

# Row 7 - Tela Suporte
$ws.Range("E7").Value = "G"
$ws.Range("F7").Value = 13

# Row 8 - Campo de pesquisa
$ws.Range("E8").Value = "M"
$ws.Range("F8").Value = 8

# Row 9 - Tela Sobre Nós (style also switches to vertically-centered)
$ws.Range("E9").Value = "GG"
$ws.Range("F9").Value = 21
$ws.Range("E9").VerticalAlignment = -4108

# Row 10 - Barra de Menu (Tela Perfil)
$ws.Range("E10").Value = "GG"
$ws.Range("F10").Value = 21
$ws.Range("E10").VerticalAlignment = -4108

# Row 11 - Tela Perfil
$ws.Range("E11").Value = "GG"
$ws.Range("F11").Value = 21
$ws.Range("E11").VerticalAlignment = -4108

# Row 12 - SIMULADOR FINANCEIRO
$ws.Range("E12").Value = "P"
$ws.Range("F12").Value = 5

# Row 13 - Tela DASHBOARD
$ws.Range("E13").Value = "P"
$ws.Range("F13").Value = 5

# Row 14 - Banco de dados MYSQL (style also switches to vertically-centered)
$ws.Range("E14").Value = "M"
$ws.Range("F14").Value = 8
$ws.Range("E14").VerticalAlignment = -4108

# Row 15 - Modelagem lógica
$ws.Range("E15").Value = "P"
$ws.Range("F15").Value = 5

# Row 16 - Script banco de dados
$ws.Range("E16").Value = "P"
$ws.Range("F16").Value = 5

# Row 17 - Planilha de riscos: only the TAMANHO cell formatting changes
# (vertical-centered) - no size value is entered for this row.
$ws.Range("E17").VerticalAlignment = -4108

# Row 18 - Diagrama de solução (style also switches to vertically-centered)
$ws.Range("E18").Value = "G"
$ws.Range("F18").Value = 13
$ws.Range("E18").VerticalAlignment = -4108

# Row 19 - Teste com os sensores
$ws.Range("E19").Value = "G"
$ws.Range("F19").Value = 13

# Restore the selection to reflect where the user left off (ATAS review).
[void]$ws.Range("G2").Select()
